$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2025
$ws.Range("B2").Value = "JAN"
$ws.Range("C2").Value = "10/01-01/01"
$ws.Range("D2").Value = "11/01, 12/01, 13/01, 14/01, 15/01, 16/01, 17/01, 18/01, 19/01, 20/01, 21/01, 22/01, 23/01, 24/01, 25/01, 26/01, 27/01, 28/01, 29/01, 30/01, 31/01"
